$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.728.12"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.027.31"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "2.325.93"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.758"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("D17").Value = "2.022.93"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "37.648.17"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -6.60%  "
$ws.Range("D21").Value = "0.0₃0821"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  +7.44%  "
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("D41").Value = "1.531.90"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0907"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "2.215.46"
$ws.Range("E51").Value = "  -1.96%  "
